$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129, shifting existing rows 129-214 down to 130-215
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with this week's price data for
# Feria Lagunitas de Puerto Montt - Ajo Chino, Primera
$ws.Cells.Item(129, 1).Value2 = 4
$ws.Cells.Item(129, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(129, 3).Value2 = "Los Lagos"
$ws.Cells.Item(129, 4).Value2 = 44596
$ws.Cells.Item(129, 5).Value2 = 10
$ws.Cells.Item(129, 6).Value2 = 100112003
$ws.Cells.Item(129, 7).Value2 = "Ajo"
$ws.Cells.Item(129, 8).Value2 = "Chino"
$ws.Cells.Item(129, 9).Value2 = "Primera"
$ws.Cells.Item(129, 10).Value2 = 180
$ws.Cells.Item(129, 11).Value2 = 20000
$ws.Cells.Item(129, 12).Value2 = 21000
$ws.Cells.Item(129, 13).Value2 = 20500
$ws.Cells.Item(129, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(129, 15).Value2 = "China"
$ws.Cells.Item(129, 16).Value2 = 2050
$ws.Cells.Item(129, 17).Value2 = 10
$ws.Cells.Item(129, 18).Value2 = "Hortaliza"
